# Updates the crypto price ("D") and 1h volume-change ("E") columns to match
# the freshly scraped figures (GitHub Actions refresh).

function Set-CellText($range, [string]$value) {
    # The "Price" column stores plain text (e.g. "1.00", "0.999", "66.105.47")
    # rather than numbers -- some of those look numeric to Excel, which would
    # otherwise auto-convert them (dropping trailing zeros, etc.) when assigned
    # through .Value. A leading apostrophe forces literal text entry; the
    # cell's style is restored right after so no stray quote-prefix formatting
    # is left on the cell.
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws.Range("D2") "66.105.47"
$ws.Range("E2").Value = "  -2.63%  "

Set-CellText $ws.Range("D3") "3.827.53"
$ws.Range("E3").Value = "  +1.57%  "

Set-CellText $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.11%  "

Set-CellText $ws.Range("D5") "421.68"
$ws.Range("E5").Value = "  +0.07%  "

Set-CellText $ws.Range("D6") "127.33"
$ws.Range("E6").Value = "  -3.90%  "

Set-CellText $ws.Range("D7") "3.819.50"
$ws.Range("E7").Value = "  +1.84%  "

Set-CellText $ws.Range("D8") "0.601"
$ws.Range("E8").Value = "  -7.91%  "

Set-CellText $ws.Range("D9") "1.00"
$ws.Range("E9").Value = "  +0.09%  "

Set-CellText $ws.Range("D10") "0.713"
$ws.Range("E10").Value = "  -8.11%  "

Set-CellText $ws.Range("D11") "0.162"
$ws.Range("E11").Value = "  -13.71%  "

Set-CellText $ws.Range("D12") "0.0000341"
$ws.Range("E12").Value = "  -21.55%  "

Set-CellText $ws.Range("D13") "39.93"
$ws.Range("E13").Value = "  -7.17%  "

Set-CellText $ws.Range("D14") "4.435.98"
$ws.Range("E14").Value = "  +1.78%  "

Set-CellText $ws.Range("D15") "9.85"
$ws.Range("E15").Value = "  -5.28%  "

Set-CellText $ws.Range("D16") "15.75"
$ws.Range("E16").Value = "  +20.45%  "

Set-CellText $ws.Range("D17") "3.843.82"
$ws.Range("E17").Value = "  +1.79%  "

$ws.Range("E18").Value = "  -1.95%  "

Set-CellText $ws.Range("D19") "19.39"
$ws.Range("E19").Value = "  -6.09%  "

Set-CellText $ws.Range("D20") "66.447.58"
$ws.Range("E20").Value = "  -2.12%  "

$ws.Range("E21").Value = "  -7.22%  "

Set-CellText $ws.Range("D22") "400.06"
$ws.Range("E22").Value = "  -11.53%  "

Set-CellText $ws.Range("D23") "14.13"
$ws.Range("E23").Value = "  -11.49%  "

Set-CellText $ws.Range("D24") "83.34"
$ws.Range("E24").Value = "  -8.35%  "

Set-CellText $ws.Range("D25") "2.95"
$ws.Range("E25").Value = "  -5.01%  "

Set-CellText $ws.Range("D26") "36.71"
$ws.Range("E26").Value = "  -4.72%  "

$ws.Range("E27").Value = "  +12.54%  "

$ws.Range("E28").Value = "  -5.42%  "

Set-CellText $ws.Range("D29") "9.32"
$ws.Range("E29").Value = "  -8.13%  "

Set-CellText $ws.Range("D30") "700.25"
$ws.Range("E30").Value = "  +2.14%  "

$ws.Range("E31").Value = "  +1.55%  "

$ws.Range("E32").Value = "  -4.18%  "

Set-CellText $ws.Range("D33") "12.19"
$ws.Range("E33").Value = "  -4.26%  "

Set-CellText $ws.Range("D34") "7.43"
$ws.Range("E34").Value = "  +3.17%  "

$ws.Range("E35").Value = "  -10.68%  "

Set-CellText $ws.Range("D36") "37.51"
$ws.Range("E36").Value = "  -11.13%  "

Set-CellText $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  +0.00%  "

Set-CellText $ws.Range("D38") "54.56"
$ws.Range("E38").Value = "  -5.47%  "

Set-CellText $ws.Range("D39") "0.0₃0751"
$ws.Range("E39").Value = "  -1.08%  "

Set-CellText $ws.Range("D40") "0.0449"
$ws.Range("E40").Value = "  -9.40%  "

Set-CellText $ws.Range("D41") "2.90"
$ws.Range("E41").Value = "  -4.35%  "

$ws.Range("E42").Value = "  +0.56%  "

Set-CellText $ws.Range("D43") "0.133"
$ws.Range("E43").Value = "  -10.48%  "

Set-CellText $ws.Range("D44") "4.38"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("E45").Value = "  -3.32%  "

Set-CellText $ws.Range("D46") "143.86"
$ws.Range("E46").Value = "  -2.64%  "

Set-CellText $ws.Range("D47") "3.09"
$ws.Range("E47").Value = "  -2.33%  "

Set-CellText $ws.Range("D48") "2.04"
$ws.Range("E48").Value = "  -4.39%  "

Set-CellText $ws.Range("D49") "25.94"
$ws.Range("E49").Value = "  -7.06%  "

$ws.Range("E50").Value = "  -5.12%  "

$ws.Range("E51").Value = "  -7.57%  "
